$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values move together per row: D (date), J, K, L, M, P
$cols = @("D", "J", "K", "L", "M", "P")

# New row -> source (old) row mapping, derived from the diff
$mapping = @{
    2  = 4
    3  = 2
    4  = 8
    5  = 10
    6  = 9
    7  = 6
    8  = 3
    9  = 5
    10 = 7
}

# Snapshot the original values for rows 2-10 before writing anything,
# since several rows both give and receive values (permutation).
# NOTE: use Value2 (not Value) for reading - Value getter misbehaves in
# this runtime and returns a descriptor string instead of the real value.
$original = @{}
foreach ($row in 2..10) {
    $original[$row] = @{}
    foreach ($col in $cols) {
        $original[$row][$col] = $ws.Range("$col$row").Value2
    }
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value2 = $original[$oldRow][$col]
    }
}
